$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column F (dSF)
$updates = @{
    2  = 6
    4  = -2
    5  = -1
    6  = -3
    7  = -2
    8  = -1
    10 = -4
    12 = 5
    13 = 3
    14 = -1
    15 = -6
    16 = 2
    18 = 5
    20 = -7
    21 = -2
    22 = 1
    23 = 1
    24 = -3
    25 = 2
    27 = 2
    29 = -2
    31 = 3
    32 = -1
    33 = 2
    34 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
